$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 234
$ws.Range("F3").Value = 1319
$ws.Range("F4").Value = 18
$ws.Range("F5").Value = 881
$ws.Range("F6").Value = 30
$ws.Range("F7").Value = 1185
$ws.Range("F8").Value = 1488
$ws.Range("F9").Value = 148
$ws.Range("F11").Value = 486
$ws.Range("F12").Value = 427
$ws.Range("F13").Value = 91
$ws.Range("F17").Value = 76
$ws.Range("F18").Value = 5874
$ws.Range("F20").Value = 5684
$ws.Range("F21").Value = 9696
$ws.Range("F23").Value = 167
$ws.Range("F24").Value = 169
$ws.Range("F25").Value = 260
$ws.Range("F26").Value = 481
$ws.Range("F27").Value = 155
$ws.Range("F28").Value = 135
$ws.Range("F29").Value = 4342
$ws.Range("F30").Value = 348

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 1158
$ws.Range("F16").Value = 92

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 608

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 608
$ws.Range("F4").Value = 234
$ws.Range("F5").Value = 1319
$ws.Range("F6").Value = 18
$ws.Range("F8").Value = 881
$ws.Range("F9").Value = 30
$ws.Range("F10").Value = 1185
$ws.Range("F12").Value = 1488
$ws.Range("F14").Value = 148
$ws.Range("F15").Value = 486
$ws.Range("F17").Value = 427
$ws.Range("F18").Value = 91
$ws.Range("F23").Value = 76
$ws.Range("F24").Value = 5874
$ws.Range("F26").Value = 5684
$ws.Range("F27").Value = 9696
$ws.Range("F30").Value = 167
$ws.Range("F31").Value = 169
$ws.Range("F32").Value = 260
$ws.Range("F34").Value = 481
$ws.Range("F37").Value = 155
$ws.Range("F38").Value = 135
$ws.Range("F39").Value = 4342
$ws.Range("F42").Value = 92
$ws.Range("F46").Value = 348
